$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0.6672711372375488
    3  = 0.4321568012237549
    4  = 0.6069910526275635
    5  = 0.7203068733215332
    6  = 0.224308967590332
    7  = 0.2178268432617188
    8  = 0.2365880012512207
    9  = 0.23577880859375
    10 = 0.3554270267486572
    11 = 0.2591707706451416
    12 = 1.180288791656494
    13 = 0.6191408634185791
    14 = 0.4552199840545654
    15 = 0.4576971530914307
    16 = 0.7763431072235107
    17 = 0.6328058242797852
    18 = 0.2371830940246582
    19 = 0.2255909442901611
    20 = 0.2602841854095459
    21 = 0.2583858966827393
    22 = 0.2893581390380859
    23 = 0.2859761714935303
    24 = 0.6115460395812988
    25 = 0.8260171413421631
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
